$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44239
$ws.Cells.Item(2, 10).Value = 600
$ws.Cells.Item(2, 11).Value = 28000
$ws.Cells.Item(2, 12).Value = 29000
$ws.Cells.Item(2, 13).Value = 28500
$ws.Cells.Item(2, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(2, 16).Value = 1140

$ws.Cells.Item(3, 4).Value = 44281
$ws.Cells.Item(3, 10).Value = 640
$ws.Cells.Item(3, 11).Value = 29000
$ws.Cells.Item(3, 12).Value = 30000
$ws.Cells.Item(3, 13).Value = 29500
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 1180

$ws.Cells.Item(4, 4).Value = 44279
$ws.Cells.Item(4, 10).Value = 560
$ws.Cells.Item(4, 11).Value = 29000
$ws.Cells.Item(4, 12).Value = 30000
$ws.Cells.Item(4, 13).Value = 29500
$ws.Cells.Item(4, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(4, 16).Value = 1180

$ws.Cells.Item(5, 4).Value = 44221
$ws.Cells.Item(5, 10).Value = 460
$ws.Cells.Item(5, 11).Value = 35000
$ws.Cells.Item(5, 12).Value = 36000
$ws.Cells.Item(5, 13).Value = 35500
$ws.Cells.Item(5, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(5, 16).Value = 1420

$ws.Cells.Item(6, 4).Value = 44251
$ws.Cells.Item(6, 10).Value = 700
$ws.Cells.Item(6, 11).Value = 29000
$ws.Cells.Item(6, 12).Value = 30000
$ws.Cells.Item(6, 13).Value = 29500
$ws.Cells.Item(6, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(6, 16).Value = 1180

$ws.Cells.Item(7, 4).Value = 44277
$ws.Cells.Item(7, 10).Value = 560
$ws.Cells.Item(7, 11).Value = 29000
$ws.Cells.Item(7, 12).Value = 30000
$ws.Cells.Item(7, 13).Value = 29500
$ws.Cells.Item(7, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(7, 16).Value = 1180

$ws.Cells.Item(8, 4).Value = 44237
$ws.Cells.Item(8, 10).Value = 600
$ws.Cells.Item(8, 11).Value = 29000
$ws.Cells.Item(8, 12).Value = 30000
$ws.Cells.Item(8, 13).Value = 29500
$ws.Cells.Item(8, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(8, 16).Value = 1180

$ws.Cells.Item(9, 4).Value = 44302
$ws.Cells.Item(9, 10).Value = 600
$ws.Cells.Item(9, 11).Value = 28000
$ws.Cells.Item(9, 12).Value = 29000
$ws.Cells.Item(9, 13).Value = 28500
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 1140

$ws.Cells.Item(10, 4).Value = 44309
$ws.Cells.Item(10, 10).Value = 600
$ws.Cells.Item(10, 11).Value = 28000
$ws.Cells.Item(10, 12).Value = 29000
$ws.Cells.Item(10, 13).Value = 28500
$ws.Cells.Item(10, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 16).Value = 1140

$ws.Cells.Item(11, 4).Value = 44253
$ws.Cells.Item(11, 10).Value = 660
$ws.Cells.Item(11, 11).Value = 28000
$ws.Cells.Item(11, 12).Value = 29000
$ws.Cells.Item(11, 13).Value = 28500
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 1140

$ws.Cells.Item(12, 4).Value = 44323
$ws.Cells.Item(12, 10).Value = 600
$ws.Cells.Item(12, 11).Value = 29000
$ws.Cells.Item(12, 12).Value = 30000
$ws.Cells.Item(12, 13).Value = 29500
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 1180

$ws.Cells.Item(13, 4).Value = 44238
$ws.Cells.Item(13, 10).Value = 520
$ws.Cells.Item(13, 11).Value = 28000
$ws.Cells.Item(13, 12).Value = 29000
$ws.Cells.Item(13, 13).Value = 28500
$ws.Cells.Item(13, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(13, 16).Value = 1140

$ws.Cells.Item(14, 4).Value = 44260
$ws.Cells.Item(14, 10).Value = 600
$ws.Cells.Item(14, 11).Value = 27000
$ws.Cells.Item(14, 12).Value = 28000
$ws.Cells.Item(14, 13).Value = 27500
$ws.Cells.Item(14, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 16).Value = 1100

$ws.Cells.Item(15, 4).Value = 44272
$ws.Cells.Item(15, 10).Value = 600
$ws.Cells.Item(15, 11).Value = 28000
$ws.Cells.Item(15, 12).Value = 29000
$ws.Cells.Item(15, 13).Value = 28500
$ws.Cells.Item(15, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(15, 16).Value = 1140

$ws.Cells.Item(16, 4).Value = 44218
$ws.Cells.Item(16, 10).Value = 400
$ws.Cells.Item(16, 11).Value = 34000
$ws.Cells.Item(16, 12).Value = 35000
$ws.Cells.Item(16, 13).Value = 34500
$ws.Cells.Item(16, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(16, 16).Value = 1380

$ws.Cells.Item(17, 4).Value = 44200
$ws.Cells.Item(17, 10).Value = 400
$ws.Cells.Item(17, 11).Value = 31000
$ws.Cells.Item(17, 12).Value = 32000
$ws.Cells.Item(17, 13).Value = 31500
$ws.Cells.Item(17, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 16).Value = 1260

$ws.Cells.Item(18, 4).Value = 44246
$ws.Cells.Item(18, 10).Value = 600
$ws.Cells.Item(18, 11).Value = 28000
$ws.Cells.Item(18, 12).Value = 29000
$ws.Cells.Item(18, 13).Value = 28500
$ws.Cells.Item(18, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 16).Value = 1140

$ws.Cells.Item(19, 4).Value = 44295
$ws.Cells.Item(19, 10).Value = 600
$ws.Cells.Item(19, 11).Value = 29000
$ws.Cells.Item(19, 12).Value = 30000
$ws.Cells.Item(19, 13).Value = 29500
$ws.Cells.Item(19, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 16).Value = 1180

$ws.Cells.Item(21, 4).Value = 44278
$ws.Cells.Item(21, 10).Value = 400
$ws.Cells.Item(21, 11).Value = 29000
$ws.Cells.Item(21, 12).Value = 30000
$ws.Cells.Item(21, 13).Value = 29500
$ws.Cells.Item(21, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(21, 16).Value = 1180

$ws.Cells.Item(22, 4).Value = 44245
$ws.Cells.Item(22, 10).Value = 540
$ws.Cells.Item(22, 11).Value = 28000
$ws.Cells.Item(22, 12).Value = 29000
$ws.Cells.Item(22, 13).Value = 28500
$ws.Cells.Item(22, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(22, 16).Value = 1140

$ws.Cells.Item(23, 4).Value = 44230
$ws.Cells.Item(23, 10).Value = 500
$ws.Cells.Item(23, 11).Value = 28000
$ws.Cells.Item(23, 12).Value = 29000
$ws.Cells.Item(23, 13).Value = 28500
$ws.Cells.Item(23, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(23, 16).Value = 1140

$ws.Cells.Item(24, 4).Value = 44216
$ws.Cells.Item(24, 10).Value = 600
$ws.Cells.Item(24, 11).Value = 36000
$ws.Cells.Item(24, 12).Value = 37000
$ws.Cells.Item(24, 13).Value = 36500
$ws.Cells.Item(24, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(24, 16).Value = 1460

$ws.Cells.Item(25, 4).Value = 44244
$ws.Cells.Item(25, 10).Value = 640
$ws.Cells.Item(25, 11).Value = 29000
$ws.Cells.Item(25, 12).Value = 30000
$ws.Cells.Item(25, 13).Value = 29500
$ws.Cells.Item(25, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(25, 16).Value = 1180

$ws.Cells.Item(26, 4).Value = 44330
$ws.Cells.Item(26, 10).Value = 520
$ws.Cells.Item(26, 11).Value = 32000
$ws.Cells.Item(26, 12).Value = 33000
$ws.Cells.Item(26, 13).Value = 32500
$ws.Cells.Item(26, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(26, 16).Value = 1300

$ws.Cells.Item(27, 4).Value = 44265
$ws.Cells.Item(27, 10).Value = 760
$ws.Cells.Item(27, 11).Value = 28000
$ws.Cells.Item(27, 12).Value = 29000
$ws.Cells.Item(27, 13).Value = 28500
$ws.Cells.Item(27, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(27, 16).Value = 1140

$ws.Cells.Item(28, 4).Value = 44264
$ws.Cells.Item(28, 10).Value = 400
$ws.Cells.Item(28, 11).Value = 28000
$ws.Cells.Item(28, 12).Value = 29000
$ws.Cells.Item(28, 13).Value = 28500
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 1140

$ws.Cells.Item(29, 4).Value = 44293
$ws.Cells.Item(29, 10).Value = 500
$ws.Cells.Item(29, 11).Value = 29000
$ws.Cells.Item(29, 12).Value = 30000
$ws.Cells.Item(29, 13).Value = 29500
$ws.Cells.Item(29, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(29, 16).Value = 1180

$ws.Cells.Item(30, 4).Value = 44209
$ws.Cells.Item(30, 10).Value = 600
$ws.Cells.Item(30, 11).Value = 36000
$ws.Cells.Item(30, 12).Value = 37000
$ws.Cells.Item(30, 13).Value = 36500
$ws.Cells.Item(30, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(30, 16).Value = 1460

$ws.Cells.Item(31, 4).Value = 44316
$ws.Cells.Item(31, 10).Value = 600
$ws.Cells.Item(31, 11).Value = 28000
$ws.Cells.Item(31, 12).Value = 29000
$ws.Cells.Item(31, 13).Value = 28500
$ws.Cells.Item(31, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(31, 16).Value = 1140

$ws.Cells.Item(32, 4).Value = 44291
$ws.Cells.Item(32, 10).Value = 500
$ws.Cells.Item(32, 11).Value = 29000
$ws.Cells.Item(32, 12).Value = 30000
$ws.Cells.Item(32, 13).Value = 29500
$ws.Cells.Item(32, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(32, 16).Value = 1180

$ws.Cells.Item(33, 4).Value = 44274
$ws.Cells.Item(33, 10).Value = 600
$ws.Cells.Item(33, 11).Value = 28000
$ws.Cells.Item(33, 12).Value = 29000
$ws.Cells.Item(33, 13).Value = 28500
$ws.Cells.Item(33, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(33, 16).Value = 1140

$ws.Cells.Item(34, 4).Value = 44335
$ws.Cells.Item(34, 10).Value = 400
$ws.Cells.Item(34, 11).Value = 30000
$ws.Cells.Item(34, 12).Value = 31000
$ws.Cells.Item(34, 13).Value = 30500
$ws.Cells.Item(34, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 16).Value = 1220

$ws.Cells.Item(35, 4).Value = 44231
$ws.Cells.Item(35, 10).Value = 500
$ws.Cells.Item(35, 11).Value = 25000
$ws.Cells.Item(35, 12).Value = 26000
$ws.Cells.Item(35, 13).Value = 25500
$ws.Cells.Item(35, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(35, 16).Value = 1020

$ws.Cells.Item(37, 4).Value = 44300
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 29000
$ws.Cells.Item(37, 12).Value = 30000
$ws.Cells.Item(37, 13).Value = 29500
$ws.Cells.Item(37, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 16).Value = 1180

$ws.Cells.Item(38, 4).Value = 44202
$ws.Cells.Item(38, 10).Value = 600
$ws.Cells.Item(38, 11).Value = 29000
$ws.Cells.Item(38, 12).Value = 30000
$ws.Cells.Item(38, 13).Value = 29500
$ws.Cells.Item(38, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(38, 16).Value = 1180

$ws.Cells.Item(39, 4).Value = 44307
$ws.Cells.Item(39, 10).Value = 560
$ws.Cells.Item(39, 11).Value = 28000
$ws.Cells.Item(39, 12).Value = 29000
$ws.Cells.Item(39, 13).Value = 28500
$ws.Cells.Item(39, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(39, 16).Value = 1140

$ws.Cells.Item(40, 4).Value = 44321
$ws.Cells.Item(40, 10).Value = 400
$ws.Cells.Item(40, 11).Value = 29000
$ws.Cells.Item(40, 12).Value = 30000
$ws.Cells.Item(40, 13).Value = 29500
$ws.Cells.Item(40, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(40, 16).Value = 1180

$ws.Cells.Item(41, 4).Value = 44223
$ws.Cells.Item(41, 10).Value = 660
$ws.Cells.Item(41, 11).Value = 32500
$ws.Cells.Item(41, 12).Value = 33000
$ws.Cells.Item(41, 13).Value = 32750
$ws.Cells.Item(41, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(41, 16).Value = 1310

$ws.Cells.Item(42, 4).Value = 44298
$ws.Cells.Item(42, 10).Value = 500
$ws.Cells.Item(42, 11).Value = 29000
$ws.Cells.Item(42, 12).Value = 30000
$ws.Cells.Item(42, 13).Value = 29500
$ws.Cells.Item(42, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 16).Value = 1180

$ws.Cells.Item(43, 4).Value = 44333
$ws.Cells.Item(43, 10).Value = 400
$ws.Cells.Item(43, 11).Value = 31500
$ws.Cells.Item(43, 12).Value = 32000
$ws.Cells.Item(43, 13).Value = 31750
$ws.Cells.Item(43, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(43, 16).Value = 1270

$ws.Cells.Item(44, 4).Value = 44187
$ws.Cells.Item(44, 10).Value = 400
$ws.Cells.Item(44, 11).Value = 37000
$ws.Cells.Item(44, 12).Value = 38000
$ws.Cells.Item(44, 13).Value = 37500
$ws.Cells.Item(44, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(44, 16).Value = 1500

$ws.Cells.Item(45, 4).Value = 44326
$ws.Cells.Item(45, 10).Value = 500
$ws.Cells.Item(45, 11).Value = 29000
$ws.Cells.Item(45, 12).Value = 30000
$ws.Cells.Item(45, 13).Value = 29500
$ws.Cells.Item(45, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(45, 16).Value = 1180
